$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(236).Insert()

$ws.Cells.Item(236, 1).Value = 8
$ws.Cells.Item(236, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44522
$ws.Cells.Item(236, 5).Value = 4
$ws.Cells.Item(236, 6).Value = 100114001
$ws.Cells.Item(236, 7).Value = "Papa"
$ws.Cells.Item(236, 8).Value = "Cardinal"
$ws.Cells.Item(236, 9).Value = "1a nueva(o)"
$ws.Cells.Item(236, 10).Value = 2000
$ws.Cells.Item(236, 11).Value = 11500
$ws.Cells.Item(236, 12).Value = 12000
$ws.Cells.Item(236, 13).Value = 11750
$ws.Cells.Item(236, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(236, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(236, 16).Value = 470
$ws.Cells.Item(236, 17).Value = 25
$ws.Cells.Item(236, 18).Value = "Hortaliza"
